$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Cells.Item(53, 8).Value = 3011.818
$ws.Cells.Item(53, 9).Value = 316.66666
$ws.Cells.Item(53, 10).Value = 4022.5
$ws.Cells.Item(53, 11).Value = 316.66666
$ws.Cells.Item(53, 12).Value = 4022.5
$ws.Cells.Item(53, 13).Value = 320.33334
$ws.Cells.Item(53, 14).Value = -5296.5
# Row 127
$ws.Cells.Item(127, 8).Value = 1512.3334
$ws.Cells.Item(127, 9).Value = 898.5
$ws.Cells.Item(127, 11).Value = 2695.5
$ws.Cells.Item(127, 13).Value = 2264.5
# Row 138
$ws.Cells.Item(138, 8).Value = 2602.0205
$ws.Cells.Item(138, 9).Value = 1971.0714
$ws.Cells.Item(138, 10).Value = 2854.4
$ws.Cells.Item(138, 11).Value = 5913.2142
$ws.Cells.Item(138, 12).Value = 8563.200000000001
$ws.Cells.Item(138, 13).Value = -773.2142000000003
$ws.Cells.Item(138, 14).Value = -18843.2
# Row 139
$ws.Cells.Item(139, 8).Value = 52780
$ws.Cells.Item(139, 10).Value = 52780
$ws.Cells.Item(139, 12).Value = 52780
$ws.Cells.Item(139, 14).Value = -63060
# Row 141
$ws.Cells.Item(141, 8).Value = 2706.842
$ws.Cells.Item(141, 9).Value = 1901.875
$ws.Cells.Item(141, 11).Value = 5705.625
$ws.Cells.Item(141, 13).Value = -525.625

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 14385.987
$ws.Cells.Item(32, 9).Value = 15482.694
$ws.Cells.Item(32, 10).Value = 5612.3335
$ws.Cells.Item(32, 11).Value = 15482.694
$ws.Cells.Item(32, 12).Value = 5612.3335
$ws.Cells.Item(32, 13).Value = -15195.694
$ws.Cells.Item(32, 14).Value = -6186.3335
# Row 97
$ws.Cells.Item(97, 8).Value = 2418.3333
$ws.Cells.Item(97, 9).Value = 2702
$ws.Cells.Item(97, 11).Value = 2702
$ws.Cells.Item(97, 13).Value = -2206
# Row 132
$ws.Cells.Item(132, 8).Value = 41885.31
$ws.Cells.Item(132, 9).Value = 3185.3333
$ws.Cells.Item(132, 10).Value = 75056.71000000001
$ws.Cells.Item(132, 11).Value = 9555.999899999999
$ws.Cells.Item(132, 12).Value = 225170.13
$ws.Cells.Item(132, 13).Value = -7025.999899999999
$ws.Cells.Item(132, 14).Value = -230230.13
# Row 135
$ws.Cells.Item(135, 8).Value = 26485.8
$ws.Cells.Item(135, 10).Value = 26485.8
$ws.Cells.Item(135, 12).Value = 26485.8
$ws.Cells.Item(135, 14).Value = -36625.8

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1395.3572
$ws.Cells.Item(20, 9).Value = 966.8182
$ws.Cells.Item(20, 11).Value = 966.8182
$ws.Cells.Item(20, 13).Value = -719.8182
# Row 64
$ws.Cells.Item(64, 8).Value = 1175.0476
$ws.Cells.Item(64, 9).Value = 2828
$ws.Cells.Item(64, 11).Value = 2828
$ws.Cells.Item(64, 13).Value = -2603
# Row 67
$ws.Cells.Item(67, 8).Value = 1175.0476
$ws.Cells.Item(67, 9).Value = 2828
$ws.Cells.Item(67, 11).Value = 2828
$ws.Cells.Item(67, 13).Value = -2048
# Row 94
$ws.Cells.Item(94, 8).Value = 2721.7896
$ws.Cells.Item(94, 9).Value = 1193.1428
$ws.Cells.Item(94, 10).Value = 7002
$ws.Cells.Item(94, 11).Value = 1193.1428
$ws.Cells.Item(94, 12).Value = 7002
$ws.Cells.Item(94, 13).Value = -742.1428000000001
$ws.Cells.Item(94, 14).Value = -7904
# Row 107
$ws.Cells.Item(107, 8).Value = 1532.5454
$ws.Cells.Item(107, 9).Value = 1538.8
$ws.Cells.Item(107, 11).Value = 1538.8
$ws.Cells.Item(107, 13).Value = 381.2

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 44.42857
$ws.Cells.Item(7, 9).Value = 43.5
$ws.Cells.Item(7, 11).Value = 43.5
$ws.Cells.Item(7, 13).Value = 69.5
# Row 22
$ws.Cells.Item(22, 8).Value = 559.1111
$ws.Cells.Item(22, 9).Value = 516.5
$ws.Cells.Item(22, 11).Value = 516.5
$ws.Cells.Item(22, 13).Value = -166.5
# Row 31
$ws.Cells.Item(31, 8).Value = 14920.576
$ws.Cells.Item(31, 9).Value = 31894.23
$ws.Cells.Item(31, 10).Value = 3887.7
$ws.Cells.Item(31, 11).Value = 31894.23
$ws.Cells.Item(31, 12).Value = 3887.7
$ws.Cells.Item(31, 13).Value = -31599.23
$ws.Cells.Item(31, 14).Value = -4477.7
# Row 34
$ws.Cells.Item(34, 8).Value = 14920.576
$ws.Cells.Item(34, 9).Value = 31894.23
$ws.Cells.Item(34, 10).Value = 3887.7
$ws.Cells.Item(34, 11).Value = 31894.23
$ws.Cells.Item(34, 12).Value = 3887.7
$ws.Cells.Item(34, 13).Value = -31692.23
$ws.Cells.Item(34, 14).Value = -4291.7
# Row 100
$ws.Cells.Item(100, 8).Value = 200000
$ws.Cells.Item(100, 10).Value = 200000
$ws.Cells.Item(100, 12).Value = 200000
$ws.Cells.Item(100, 14).Value = -202164

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 490.23215
$ws.Cells.Item(5, 9).Value = 398.03845
$ws.Cells.Item(5, 10).Value = 570.13336
$ws.Cells.Item(5, 11).Value = 1194.11535
$ws.Cells.Item(5, 12).Value = 1710.40008
$ws.Cells.Item(5, 13).Value = -1082.11535
$ws.Cells.Item(5, 14).Value = -1934.40008
# Row 12
$ws.Cells.Item(12, 8).Value = 56.615383
$ws.Cells.Item(12, 9).Value = 50
$ws.Cells.Item(12, 11).Value = 150
$ws.Cells.Item(12, 13).Value = 23
# Row 68
$ws.Cells.Item(68, 8).Value = 7379.0625
$ws.Cells.Item(68, 10).Value = 11321.6
$ws.Cells.Item(68, 12).Value = 33964.8
$ws.Cells.Item(68, 14).Value = -35586.8
# Row 71
$ws.Cells.Item(71, 8).Value = 7379.0625
$ws.Cells.Item(71, 10).Value = 11321.6
$ws.Cells.Item(71, 12).Value = 101894.4
$ws.Cells.Item(71, 14).Value = -110006.4
# Row 107
$ws.Cells.Item(107, 8).Value = 4439.4644
$ws.Cells.Item(107, 10).Value = 969.8261
$ws.Cells.Item(107, 12).Value = 2909.4783
$ws.Cells.Item(107, 14).Value = -6749.4783
# Row 122
$ws.Cells.Item(122, 8).Value = 1066.5385
$ws.Cells.Item(122, 10).Value = 1606.7333
$ws.Cells.Item(122, 12).Value = 14460.5997
$ws.Cells.Item(122, 14).Value = -19360.5997
# Row 131
$ws.Cells.Item(131, 8).Value = 121299.31
$ws.Cells.Item(131, 9).Value = 693.75
$ws.Cells.Item(131, 10).Value = 134163.9
$ws.Cells.Item(131, 11).Value = 2081.25
$ws.Cells.Item(131, 12).Value = 402491.7
$ws.Cells.Item(131, 13).Value = 2958.75
$ws.Cells.Item(131, 14).Value = -412571.7
# Row 135
$ws.Cells.Item(135, 8).Value = 490.23215
$ws.Cells.Item(135, 9).Value = 398.03845
$ws.Cells.Item(135, 10).Value = 570.13336
$ws.Cells.Item(135, 11).Value = 3582.34605
$ws.Cells.Item(135, 12).Value = 5131.20024
$ws.Cells.Item(135, 13).Value = -1047.34605
$ws.Cells.Item(135, 14).Value = -10201.20024

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 4884.5835
$ws.Cells.Item(7, 9).Value = 4733.727
$ws.Cells.Item(7, 11).Value = 4733.727
$ws.Cells.Item(7, 13).Value = -4621.727
# Row 55
$ws.Cells.Item(55, 8).Value = 285.66666
$ws.Cells.Item(55, 9).Value = 188.88889
$ws.Cells.Item(55, 10).Value = 382.44446
$ws.Cells.Item(55, 11).Value = 188.88889
$ws.Cells.Item(55, 12).Value = 382.44446
$ws.Cells.Item(55, 13).Value = -15.88889
$ws.Cells.Item(55, 14).Value = -728.4444599999999
# Row 93
$ws.Cells.Item(93, 8).Value = 3267.5715
$ws.Cells.Item(93, 9).Value = 3267.5715
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 3267.5715
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = -2019.5715
$ws.Cells.Item(93, 14).ClearContents()
# Row 122
$ws.Cells.Item(122, 8).Value = 3951.0833
$ws.Cells.Item(122, 9).Value = 3550.75
$ws.Cells.Item(122, 11).Value = 10652.25
$ws.Cells.Item(122, 13).Value = -8202.25
# Row 126
$ws.Cells.Item(126, 8).Value = 4884.5835
$ws.Cells.Item(126, 9).Value = 4733.727
$ws.Cells.Item(126, 11).Value = 14201.181
$ws.Cells.Item(126, 13).Value = -11731.181
# Row 132
$ws.Cells.Item(132, 8).Value = 2618.7058
$ws.Cells.Item(132, 9).Value = 1722
$ws.Cells.Item(132, 10).Value = 2992.3333
$ws.Cells.Item(132, 11).Value = 5166
$ws.Cells.Item(132, 12).Value = 8976.999899999999
$ws.Cells.Item(132, 13).Value = -2636
$ws.Cells.Item(132, 14).Value = -14036.9999
# Row 136
$ws.Cells.Item(136, 8).Value = 26238.143
$ws.Cells.Item(136, 9).Value = 47151.547
$ws.Cells.Item(136, 10).Value = 3233.4
$ws.Cells.Item(136, 11).Value = 141454.641
$ws.Cells.Item(136, 12).Value = 9700.200000000001
$ws.Cells.Item(136, 13).Value = -138904.641
$ws.Cells.Item(136, 14).Value = -14800.2

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells.Item(122, 8).Value = 1619.3636
$ws.Cells.Item(122, 9).Value = 1413.4445
$ws.Cells.Item(122, 11).Value = 4240.333500000001
$ws.Cells.Item(122, 13).Value = -1790.333500000001
# Row 132
$ws.Cells.Item(132, 8).Value = 2501.32
$ws.Cells.Item(132, 9).Value = 2492.3333
$ws.Cells.Item(132, 10).Value = 2509.6155
$ws.Cells.Item(132, 11).Value = 7476.999899999999
$ws.Cells.Item(132, 12).Value = 7528.8465
$ws.Cells.Item(132, 13).Value = -4946.999899999999
$ws.Cells.Item(132, 14).Value = -12588.8465
# Row 136
$ws.Cells.Item(136, 8).Value = 2235.2727
$ws.Cells.Item(136, 9).Value = 1497.25
$ws.Cells.Item(136, 10).Value = 2657
$ws.Cells.Item(136, 11).Value = 4491.75
$ws.Cells.Item(136, 13).Value = -1941.75
$ws.Cells.Item(136, 14).Value = -13071
